# Auto-generated edit script updating computed market-price columns (H-N)
# across multiple worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 43
$ws.Range("H43").Value = 915.5454999999999
$ws.Range("J43").Value = 1053.5
$ws.Range("L43").Value = 1053.5
$ws.Range("N43").Value = -1191.5

$ws = $wb.Worksheets.Item("ALC")  # row 98
$ws.Range("H98").Value = 1822.5883
$ws.Range("I98").Value = 1584.5714
$ws.Range("J98").Value = 2933.3333
$ws.Range("K98").Value = 1584.5714
$ws.Range("L98").Value = 2933.3333
$ws.Range("M98").Value = -86.57140000000004
$ws.Range("N98").Value = -5929.3333

$ws = $wb.Worksheets.Item("ALC")  # row 122
$ws.Range("H122").Value = 1822.5883
$ws.Range("I122").Value = 1584.5714
$ws.Range("J122").Value = 2933.3333
$ws.Range("K122").Value = 4753.7142
$ws.Range("L122").Value = 8799.999899999999
$ws.Range("M122").Value = -2303.7142
$ws.Range("N122").Value = -13699.9999

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 712224
$ws.Range("I137").Value = 3787.5
$ws.Range("J137").Value = 1077868.6
$ws.Range("K137").Value = 11362.5
$ws.Range("L137").Value = 3233605.8
$ws.Range("M137").Value = -8812.5
$ws.Range("N137").Value = -3238705.8

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 4663.5283
$ws.Range("I138").Value = 2182
$ws.Range("J138").Value = 5643.079
$ws.Range("K138").Value = 6546
$ws.Range("L138").Value = 16929.237
$ws.Range("M138").Value = -1406
$ws.Range("N138").Value = -27209.237

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 19220.262
$ws.Range("I32").Value = 19696.69
$ws.Range("K32").Value = 19696.69
$ws.Range("M32").Value = -19409.69

$ws = $wb.Worksheets.Item("ARM")  # row 110
$ws.Range("H110").Value = 1640.1333
$ws.Range("I110").Value = 1640.1333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1640.1333
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = 404.8667
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("ARM")  # row 123
$ws.Range("H123").Value = 74347.5
$ws.Range("J123").Value = 74347.5
$ws.Range("L123").Value = 74347.5
$ws.Range("N123").Value = -84147.5

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 2049.4426
$ws.Range("I132").Value = 1893.8667
$ws.Range("J132").Value = 2487
$ws.Range("K132").Value = 5681.6001
$ws.Range("L132").Value = 7461
$ws.Range("M132").Value = -3151.6001
$ws.Range("N132").Value = -12521

$ws = $wb.Worksheets.Item("BSM")  # row 94
$ws.Range("H94").Value = 2056.182
$ws.Range("I94").Value = 2077.25
$ws.Range("K94").Value = 2077.25
$ws.Range("M94").Value = -1626.25

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 3079.6155
$ws.Range("I134").Value = 3079.6155
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9238.8465
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -6703.8465
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 31
$ws.Range("H31").Value = 644104.5600000001
$ws.Range("I31").Value = 15106.143
$ws.Range("J31").Value = 844240.4399999999
$ws.Range("K31").Value = 15106.143
$ws.Range("L31").Value = 844240.4399999999
$ws.Range("M31").Value = -14811.143
$ws.Range("N31").Value = -844830.4399999999

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 644104.5600000001
$ws.Range("I34").Value = 15106.143
$ws.Range("J34").Value = 844240.4399999999
$ws.Range("K34").Value = 15106.143
$ws.Range("L34").Value = 844240.4399999999
$ws.Range("M34").Value = -14904.143
$ws.Range("N34").Value = -844644.4399999999

$ws = $wb.Worksheets.Item("CRP")  # row 132
$ws.Range("H132").Value = 2562.25
$ws.Range("I132").Value = 2144.7673
$ws.Range("J132").Value = 4556.8887
$ws.Range("K132").Value = 6434.3019
$ws.Range("L132").Value = 13670.6661
$ws.Range("M132").Value = -3904.3019
$ws.Range("N132").Value = -18730.6661

$ws = $wb.Worksheets.Item("CRP")  # row 134
$ws.Range("H134").Value = 2138.6274
$ws.Range("I134").Value = 1650.2433
$ws.Range("J134").Value = 3429.3572
$ws.Range("K134").Value = 4950.7299
$ws.Range("L134").Value = 10288.0716
$ws.Range("M134").Value = -2415.7299
$ws.Range("N134").Value = -15358.0716

$ws = $wb.Worksheets.Item("CRP")  # row 138
$ws.Range("H138").Value = 71740
$ws.Range("J138").Value = 71740
$ws.Range("L138").Value = 71740
$ws.Range("N138").Value = -82020

$ws = $wb.Worksheets.Item("CUL")  # row 22
$ws.Range("H22").Value = 142858430
$ws.Range("I22").Value = 250000750
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 750002250
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -750002081
$ws.Range("N22").Value = -6338

$ws = $wb.Worksheets.Item("CUL")  # row 27
$ws.Range("H27").Value = 142858430
$ws.Range("I27").Value = 250000750
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 750002250
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -750002148
$ws.Range("N27").Value = -6204

$ws = $wb.Worksheets.Item("CUL")  # row 68
$ws.Range("H68").Value = 2716.3604
$ws.Range("I68").Value = 1465.3914
$ws.Range("J68").Value = 4154.975
$ws.Range("K68").Value = 4396.174199999999
$ws.Range("L68").Value = 12464.925
$ws.Range("M68").Value = -3585.174199999999
$ws.Range("N68").Value = -14086.925

$ws = $wb.Worksheets.Item("CUL")  # row 71
$ws.Range("H71").Value = 2716.3604
$ws.Range("I71").Value = 1465.3914
$ws.Range("J71").Value = 4154.975
$ws.Range("K71").Value = 13188.5226
$ws.Range("L71").Value = 37394.775
$ws.Range("M71").Value = -9132.5226
$ws.Range("N71").Value = -45506.775

$ws = $wb.Worksheets.Item("CUL")  # row 80
$ws.Range("H80").Value = 2900
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = 8700
$ws.Range("N80").Value = -10572
$ws.Range("L80").ClearContents()

$ws = $wb.Worksheets.Item("CUL")  # row 83
$ws.Range("H83").Value = 2900
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = 26100
$ws.Range("N83").Value = -35460
$ws.Range("L83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")  # row 107
$ws.Range("H107").Value = 1019.5955
$ws.Range("I107").Value = 537.4722
$ws.Range("J107").Value = 1347.0754
$ws.Range("K107").Value = 1612.4166
$ws.Range("L107").Value = 4041.2262
$ws.Range("M107").Value = 307.5834
$ws.Range("N107").Value = -7881.2262

$ws = $wb.Worksheets.Item("CUL")  # row 113
$ws.Range("H113").Value = 505.99
$ws.Range("I113").Value = 600.1842
$ws.Range("J113").Value = 448.25806
$ws.Range("K113").Value = 1800.5526
$ws.Range("L113").Value = 1344.77418
$ws.Range("M113").Value = 369.4474
$ws.Range("N113").Value = -5684.77418

$ws = $wb.Worksheets.Item("CUL")  # row 141
$ws.Range("H141").Value = 3865.0908
$ws.Range("I141").Value = 2401.5386
$ws.Range("K141").Value = 7204.6158
$ws.Range("M141").Value = -2024.6158

$ws = $wb.Worksheets.Item("GSM")  # row 132
$ws.Range("H132").Value = 34834.16
$ws.Range("I132").Value = 47193.816
$ws.Range("J132").Value = 4621.6665
$ws.Range("K132").Value = 141581.448
$ws.Range("L132").Value = 13864.9995
$ws.Range("M132").Value = -139051.448
$ws.Range("N132").Value = -18924.9995

$ws = $wb.Worksheets.Item("LTW")  # row 10
$ws.Range("H10").Value = 4501
$ws.Range("I10").Value = 1002
$ws.Range("J10").Value = 8000
$ws.Range("K10").Value = 1002
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = -862
$ws.Range("N10").Value = -8280

$ws = $wb.Worksheets.Item("LTW")  # row 68
$ws.Range("H68").Value = 3861.32
$ws.Range("I68").Value = 2722.7273
$ws.Range("J68").Value = 4755.9287
$ws.Range("K68").Value = 2722.7273
$ws.Range("L68").Value = 4755.9287
$ws.Range("M68").Value = -1973.7273
$ws.Range("N68").Value = -6253.9287

$ws = $wb.Worksheets.Item("LTW")  # row 71
$ws.Range("H71").Value = 3861.32
$ws.Range("I71").Value = 2722.7273
$ws.Range("J71").Value = 4755.9287
$ws.Range("K71").Value = 13613.6365
$ws.Range("L71").Value = 23779.6435
$ws.Range("M71").Value = -9869.636500000001
$ws.Range("N71").Value = -31267.6435

$ws = $wb.Worksheets.Item("LTW")  # row 93
$ws.Range("H93").Value = 1000
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Range("H122").Value = 5759.2046
$ws.Range("I122").Value = 5705.5557
$ws.Range("J122").Value = 6000.625
$ws.Range("K122").Value = 17116.6671
$ws.Range("L122").Value = 18001.875
$ws.Range("M122").Value = -14666.6671
$ws.Range("N122").Value = -22901.875

$ws = $wb.Worksheets.Item("LTW")  # row 136
$ws.Range("H136").Value = 2705.2415
$ws.Range("I136").Value = 1776.8975
$ws.Range("J136").Value = 4610.7896
$ws.Range("K136").Value = 5330.6925
$ws.Range("L136").Value = 13832.3688
$ws.Range("M136").Value = -2780.6925
$ws.Range("N136").Value = -18932.3688

$ws = $wb.Worksheets.Item("WVR")  # row 123
$ws.Range("H123").Value = 49001.855
$ws.Range("J123").Value = 49001.855
$ws.Range("L123").Value = 49001.855
$ws.Range("N123").Value = -58801.855

$ws = $wb.Worksheets.Item("WVR")  # row 132
$ws.Range("H132").Value = 2187.85
$ws.Range("I132").Value = 2058.6897
$ws.Range("J132").Value = 2528.3635
$ws.Range("K132").Value = 6176.0691
$ws.Range("L132").Value = 7585.0905
$ws.Range("M132").Value = -3646.0691
$ws.Range("N132").Value = -12645.0905

$ws = $wb.Worksheets.Item("WVR")  # row 136
$ws.Range("H136").Value = 7095.963
$ws.Range("I136").Value = 7594.6924
$ws.Range("J136").Value = 6632.857
$ws.Range("K136").Value = 22784.0772
$ws.Range("L136").Value = 19898.571
$ws.Range("M136").Value = -20234.0772
$ws.Range("N136").Value = -24998.571
